# "Generate Report for Archive"
# - Flip the handoff status from "Ready for handoff" to "In Translation"
#   on the Overview sheet (per-locale columns) and on each locale sheet.
# - Narrow the Status-related column(s) to match the refreshed report layout.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "zh-cn" status is column E, "de-de" status is column F (row 2 = data row)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Locale sheets: "Status" is column C (row 2 = data row)
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Narrower Status columns to fit the now-shorter text
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
